# fix: checking time values
# H15/H16/H17 held elapsed-time calculations that had drifted past 24h
# (e.g. 28:00:00). Correct them to the intended wrapped time-of-day values.
# H19's cell format is reset to General (it was carrying a redundant,
# never-applied numeric format).
# Selection is moved to H16, the cell that was at the heart of the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H15").Value = 0.041666666666666664
$ws.Range("H16").Value = 0.95833333333333337
$ws.Range("H17").Value = 0

$ws.Range("H19").ClearFormats()

[void]$ws.Range("H16").Select()
